$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as plain text (e.g. "47.15"), matching
# the source site's inline-string cells. A handful of the updated prices look
# like plain decimals ("47.16", "1.20", ...); left alone, Excel's automatic
# type detection would silently convert those into numbers (dropping
# significant trailing zeros and introducing floating-point noise), so those
# specific cells are pinned to the Text format before the new value is written.
$textForcedCells = @('D5', 'D8', 'D9', 'D11', 'D14', 'D15', 'D17', 'D19', 'D21', 'D22', 'D26', 'D27', 'D28', 'D29', 'D30', 'D32', 'D33', 'D35', 'D36', 'D41', 'D43', 'D45', 'D47')
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '36.732.72'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '2.126.66'
$ws.Range('E3').Value = '  +10.78%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '256.31'
$ws.Range('E5').Value = '  +2.86%  '
$ws.Range('E6').Value = '  -4.13%  '
$ws.Range('D8').Value = '47.16'
$ws.Range('E8').Value = '  +6.20%  '
$ws.Range('D9').Value = '59.84'
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('D11').Value = '0.0747'
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').Value = '2.437.18'
$ws.Range('E13').Value = '  +10.85%  '
$ws.Range('D14').Value = '14.36'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').Value = '0.841'
$ws.Range('E15').Value = '  +5.17%  '
$ws.Range('D16').Value = '2.130.09'
$ws.Range('E16').Value = '  +11.05%  '
$ws.Range('D17').Value = '5.14'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '36.741.92'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').Value = '73.81'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D21').Value = '13.32'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '241.90'
$ws.Range('E22').Value = '  -3.92%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  -7.33%  '
$ws.Range('D26').Value = '171.95'
$ws.Range('E26').Value = '  +2.36%  '
$ws.Range('D27').Value = '21.91'
$ws.Range('E27').Value = '  +16.62%  '
$ws.Range('D28').Value = '9.24'
$ws.Range('E28').Value = '  +4.65%  '
$ws.Range('D29').Value = '2.04'
$ws.Range('E29').Value = '  -7.13%  '
$ws.Range('D30').Value = '29.27'
$ws.Range('E30').Value = '  +64.23%  '
$ws.Range('E31').Value = '  -4.39%  '
$ws.Range('D32').Value = '4.52'
$ws.Range('E32').Value = '  -0.79%  '
$ws.Range('D33').Value = '0.0964'
$ws.Range('E33').Value = '  +13.41%  '
$ws.Range('E34').Value = '  -2.56%  '
$ws.Range('D35').Value = '2.39'
$ws.Range('E35').Value = '  +17.99%  '
$ws.Range('D36').Value = '0.953'
$ws.Range('E36').Value = '  +8.65%  '
$ws.Range('E37').Value = '  -4.72%  '
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('E40').Value = '  -9.62%  '
$ws.Range('D41').Value = '1.20'
$ws.Range('E41').Value = '  +8.32%  '
$ws.Range('E42').Value = '  -1.13%  '
$ws.Range('D43').Value = '99.23'
$ws.Range('E44').Value = '  +10.49%  '
$ws.Range('D45').Value = '16.23'
$ws.Range('E45').Value = '  -5.60%  '
$ws.Range('D46').Value = '1.359.97'
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D47').Value = '7.27'
$ws.Range('E47').Value = '  +12.74%  '
$ws.Range('E48').Value = '  +3.52%  '
$ws.Range('D49').Value = '2.329.61'
$ws.Range('E49').Value = '  +10.89%  '
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('E51').Value = '  -3.25%  '
